$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Expected_ui_elements" column (C) is being moved to the end of the
# table (after column I), with columns D:I shifting left to fill C:H.
# Cutting column C and inserting it after the (then-shifted) last column
# reproduces that rearrangement along with the column widths.
$ws.Columns("C:C").Cut() | Out-Null
$ws.Columns("J:J").Insert() | Out-Null

# The UI-element description text itself was reworded (Line of Therapy ->
# Population filter 2) and now only occupies rows 2 and 3 (row 4's copy is
# removed).
$ws.Range("I2").Value = "Manage Population filter 2"
$ws.Range("I3").Value = "You can view all, create new and edit or delete existing Population filter 2 from here"
$ws.Range("I4").ClearContents() | Out-Null

# Restore the sheet view (selection moved, frozen/scrolled column reset).
$ws.Range("D8").Select() | Out-Null
